$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 675.6
$ws.Range("I6").Value = 675.6
$ws.Range("K6").Value = 2026.8
$ws.Range("M6").Value = -1914.8
$ws.Range("H33").Value = 28571662
$ws.Range("I33").Value = 29412004
$ws.Range("K33").Value = 29412004
$ws.Range("M33").Value = -29411775
$ws.Range("H116").Value = 6344.4443
$ws.Range("I116").Value = 6300
$ws.Range("J116").Value = 6414.2856
$ws.Range("K116").Value = 6300
$ws.Range("L116").Value = 6414.2856
$ws.Range("M116").Value = -2858
$ws.Range("N116").Value = -13298.2856
$ws.Range("H129").Value = 1248.7245
$ws.Range("I129").Value = 376.07144
$ws.Range("K129").Value = 1128.21432
$ws.Range("M129").Value = 3871.78568
$ws.Range("H131").Value = 2053.611
$ws.Range("I131").Value = 547
$ws.Range("J131").Value = 3936.875
$ws.Range("K131").Value = 1641
$ws.Range("L131").Value = 11810.625
$ws.Range("M131").Value = 3399
$ws.Range("N131").Value = -21890.625
$ws.Range("H138").Value = 4809939
$ws.Range("I138").Value = 1139.1892
$ws.Range("J138").Value = 16671644
$ws.Range("K138").Value = 3417.5676
$ws.Range("L138").Value = 50014932
$ws.Range("M138").Value = 1722.4324
$ws.Range("N138").Value = -50025212

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1524.9333
$ws.Range("I2").Value = 1465.8572
$ws.Range("J2").Value = 1576.625
$ws.Range("K2").Value = 1465.8572
$ws.Range("L2").Value = 1576.625
$ws.Range("M2").Value = -1352.8572
$ws.Range("N2").Value = -1802.625
$ws.Range("H4").Value = 197.33333
$ws.Range("I4").Value = 196
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 196
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -80
$ws.Range("N4").Value = -432
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 9265.681
$ws.Range("I32").Value = 10253.5
$ws.Range("K32").Value = 10253.5
$ws.Range("M32").Value = -9966.5
$ws.Range("H63").Value = 71430940
$ws.Range("I63").Value = 71430940
$ws.Range("K63").Value = 71430940
$ws.Range("M63").Value = -71430254
$ws.Range("H66").Value = 71430940
$ws.Range("I66").Value = 71430940
$ws.Range("K66").Value = 357154700
$ws.Range("M66").Value = -357151268
$ws.Range("H116").Value = 1524.9333
$ws.Range("I116").Value = 1465.8572
$ws.Range("J116").Value = 1576.625
$ws.Range("K116").Value = 1465.8572
$ws.Range("L116").Value = 1576.625
$ws.Range("M116").Value = 828.1428000000001
$ws.Range("N116").Value = -6164.625
$ws.Range("H122").Value = 7466.65
$ws.Range("I122").Value = 8330.529
$ws.Range("K122").Value = 24991.587
$ws.Range("M122").Value = -22541.587
$ws.Range("H132").Value = 15628984
$ws.Range("I132").Value = 25003624
$ws.Range("K132").Value = 75010872
$ws.Range("M132").Value = -75008342

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1524.9333
$ws.Range("I3").Value = 1465.8572
$ws.Range("J3").Value = 1576.625
$ws.Range("K3").Value = 1465.8572
$ws.Range("L3").Value = 1576.625
$ws.Range("M3").Value = -1351.8572
$ws.Range("N3").Value = -1804.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 71787810
$ws.Range("I6").Value = 83751944
$ws.Range("K6").Value = 83751944
$ws.Range("M6").Value = -83751831
$ws.Range("H132").Value = 3640.3809
$ws.Range("I132").Value = 2964.75
$ws.Range("J132").Value = 5802.4
$ws.Range("K132").Value = 8894.25
$ws.Range("L132").Value = 17407.2
$ws.Range("M132").Value = -6364.25
$ws.Range("N132").Value = -22467.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 5849.5
$ws.Range("J49").Value = 5849.5
$ws.Range("L49").Value = 17548.5
$ws.Range("N49").Value = -17860.5
$ws.Range("H81").Value = 933.8333
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 9000
$ws.Range("N81").Value = -11246
$ws.Range("H84").Value = 933.8333
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 27000
$ws.Range("N84").Value = -38232
$ws.Range("H86").Value = 2285.7144
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2285.7144
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6857.1432
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -9229.143199999999
$ws.Range("H87").Value = 16325
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496
$ws.Range("H89").Value = 2285.7144
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2285.7144
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 20571.4296
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -32427.4296
$ws.Range("H90").Value = 16325
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480
$ws.Range("H125").Value = 4633.3335
$ws.Range("J125").Value = 5087.5
$ws.Range("L125").Value = 15262.5
$ws.Range("N125").Value = -25102.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 975
$ws.Range("I22").Value = 1020
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 1020
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -725
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 975
$ws.Range("I27").Value = 1020
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 1020
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -913
$ws.Range("N27").Value = -1114
$ws.Range("H40").Value = 5388.963
$ws.Range("I40").Value = 8876
$ws.Range("J40").Value = 3920.7368
$ws.Range("K40").Value = 8876
$ws.Range("L40").Value = 3920.7368
$ws.Range("M40").Value = -8740
$ws.Range("N40").Value = -4192.736800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 29500
$ws.Range("J54").Value = 29500
$ws.Range("L54").Value = 29500
$ws.Range("N54").Value = -30540
$ws.Range("H123").Value = 44450
$ws.Range("J123").Value = 44450
$ws.Range("L123").Value = 44450
$ws.Range("N123").Value = -54250
